$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple value corrections (covid_deaths count only)
$ws.Cells.Item(133, 3).Value = 7
$ws.Cells.Item(143, 3).Value = 6
$ws.Cells.Item(216, 3).Value = 12
$ws.Cells.Item(220, 3).Value = 9
$ws.Cells.Item(1168, 3).Value = 3

# Re-sequence rows 1547-1567 (existing rows whose date/agegrp/value shifted)
$ws.Cells.Item(1547, 1).Value = 44274
$ws.Cells.Item(1547, 2).Value = "50-59"
$ws.Cells.Item(1547, 3).Value = 1

$ws.Cells.Item(1548, 1).Value = 44274
$ws.Cells.Item(1548, 2).Value = "60-69"
$ws.Cells.Item(1548, 3).Value = 4

$ws.Cells.Item(1549, 1).Value = 44275
$ws.Cells.Item(1549, 2).Value = "60-69"
$ws.Cells.Item(1549, 3).Value = 1

$ws.Cells.Item(1550, 1).Value = 44276
$ws.Cells.Item(1550, 2).Value = "60-69"
$ws.Cells.Item(1550, 3).Value = 3

$ws.Cells.Item(1551, 1).Value = 44276
$ws.Cells.Item(1551, 2).Value = "70-79"
$ws.Cells.Item(1551, 3).Value = 3

$ws.Cells.Item(1552, 1).Value = 44276
$ws.Cells.Item(1552, 2).Value = "80+"
$ws.Cells.Item(1552, 3).Value = 3

$ws.Cells.Item(1553, 1).Value = 44277
$ws.Cells.Item(1553, 2).Value = "50-59"
$ws.Cells.Item(1553, 3).Value = 1

$ws.Cells.Item(1554, 1).Value = 44277
$ws.Cells.Item(1554, 2).Value = "60-69"
$ws.Cells.Item(1554, 3).Value = 1

$ws.Cells.Item(1555, 1).Value = 44277
$ws.Cells.Item(1555, 2).Value = "70-79"
$ws.Cells.Item(1555, 3).Value = 2

$ws.Cells.Item(1556, 1).Value = 44277
$ws.Cells.Item(1556, 2).Value = "80+"
$ws.Cells.Item(1556, 3).Value = 2

$ws.Cells.Item(1557, 1).Value = 44278
$ws.Cells.Item(1557, 2).Value = "40-49"
$ws.Cells.Item(1557, 3).Value = 1

$ws.Cells.Item(1558, 1).Value = 44278
$ws.Cells.Item(1558, 2).Value = "60-69"
$ws.Cells.Item(1558, 3).Value = 1

$ws.Cells.Item(1559, 1).Value = 44278
$ws.Cells.Item(1559, 2).Value = "70-79"
$ws.Cells.Item(1559, 3).Value = 2

$ws.Cells.Item(1560, 1).Value = 44278
$ws.Cells.Item(1560, 2).Value = "80+"
$ws.Cells.Item(1560, 3).Value = 1

$ws.Cells.Item(1561, 1).Value = 44279
$ws.Cells.Item(1561, 2).Value = "60-69"
$ws.Cells.Item(1561, 3).Value = 1

$ws.Cells.Item(1562, 1).Value = 44279
$ws.Cells.Item(1562, 2).Value = "70-79"
$ws.Cells.Item(1562, 3).Value = 1

$ws.Cells.Item(1563, 1).Value = 44279
$ws.Cells.Item(1563, 2).Value = "80+"
$ws.Cells.Item(1563, 3).Value = 2

$ws.Cells.Item(1564, 1).Value = 44280
$ws.Cells.Item(1564, 2).Value = "70-79"
$ws.Cells.Item(1564, 3).Value = 5

$ws.Cells.Item(1565, 1).Value = 44280
$ws.Cells.Item(1565, 2).Value = "80+"
$ws.Cells.Item(1565, 3).Value = 2

$ws.Cells.Item(1566, 1).Value = 44281
$ws.Cells.Item(1566, 2).Value = "30-39"
$ws.Cells.Item(1566, 3).Value = 1

$ws.Cells.Item(1567, 1).Value = 44281
$ws.Cells.Item(1567, 2).Value = "50-59"
$ws.Cells.Item(1567, 3).Value = 1

# Append new rows 1568-1578
$ws.Cells.Item(1568, 1).Value = 44281
$ws.Cells.Item(1568, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1568, 2).Value = "60-69"
$ws.Cells.Item(1568, 3).Value = 1

$ws.Cells.Item(1569, 1).Value = 44281
$ws.Cells.Item(1569, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1569, 2).Value = "80+"
$ws.Cells.Item(1569, 3).Value = 1

$ws.Cells.Item(1570, 1).Value = 44282
$ws.Cells.Item(1570, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1570, 2).Value = "60-69"
$ws.Cells.Item(1570, 3).Value = 1

$ws.Cells.Item(1571, 1).Value = 44282
$ws.Cells.Item(1571, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1571, 2).Value = "70-79"
$ws.Cells.Item(1571, 3).Value = 1

$ws.Cells.Item(1572, 1).Value = 44282
$ws.Cells.Item(1572, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1572, 2).Value = "80+"
$ws.Cells.Item(1572, 3).Value = 3

$ws.Cells.Item(1573, 1).Value = 44283
$ws.Cells.Item(1573, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1573, 2).Value = "40-49"
$ws.Cells.Item(1573, 3).Value = 1

$ws.Cells.Item(1574, 1).Value = 44283
$ws.Cells.Item(1574, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1574, 2).Value = "50-59"
$ws.Cells.Item(1574, 3).Value = 1

$ws.Cells.Item(1575, 1).Value = 44283
$ws.Cells.Item(1575, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1575, 2).Value = "60-69"
$ws.Cells.Item(1575, 3).Value = 1

$ws.Cells.Item(1576, 1).Value = 44283
$ws.Cells.Item(1576, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1576, 2).Value = "70-79"
$ws.Cells.Item(1576, 3).Value = 1

$ws.Cells.Item(1577, 1).Value = 44284
$ws.Cells.Item(1577, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1577, 2).Value = "50-59"
$ws.Cells.Item(1577, 3).Value = 1

$ws.Cells.Item(1578, 1).Value = 44284
$ws.Cells.Item(1578, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1578, 2).Value = "60-69"
$ws.Cells.Item(1578, 3).Value = 1

